# Apply the "overview of data and data processing issue #2" edit:
#  - add a running-index column O (1..16) alongside the existing P:AF block
#    for the already-present rows 2-10
#  - tweak a few existing values in rows 9-10 (re-run results)
#  - append 7 brand-new data rows (11-17), each a full P:AF block plus its
#    own O index, extending the table
#  - move the active selection to AC30 (and drop the stale frozen
#    top-left-cell scroll position that pointed at Q1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New column O: sequential index 1..16 for rows 2..17 ------------------
for ($r = 2; $r -le 17; $r++) {
    $ws.Range("O$r").Value = $r - 1
}

# --- Small value corrections on existing rows ------------------------------
$ws.Range("AD9").Value = 8

$ws.Range("Q10").Value = 3
$ws.Range("Y10").Value = $false
$ws.Range("AC10").Value = 3
$ws.Range("AF10").Value = 0.23479000910422301

# --- Brand-new rows 11-17 (columns O through AF) ---------------------------
# Columns, in order: O P Q R S T U V W X Y Z AA AB AC AD AE AF
$newRows = @(
    @(10,100,3,12,1,1,0,0,1,100,$true, 0.1,6,200,3,7,3000000,0.229657338136141),
    @(11,100,3,12,1,1,0,0,1,100,$false,0.1,6,200,3,8,3000000,0.23554869440156201),
    @(12,100,3,12,1,1,0,0,1,100,$true, 0.1,6,200,3,8,3000000,0.23066012218131901),
    @(13,100,3,12,1,1,0,0,1,100,$false,0.1,6,200,5,7,3000000,0.23907163308659499),
    @(14,100,3,12,1,1,0,0,1,100,$true, 0.1,6,200,5,7,3000000,0.23371465516103901),
    @(15,100,3,12,1,1,0,0,1,100,$false,0.1,6,200,5,8,3000000,0.239124411194236),
    @(16,100,3,12,1,1,0,0,1,100,$true, 0.1,6,200,5,8,3000000,0.23333201388064201)
)

$rowIdx = 11
foreach ($rowValues in $newRows) {
    $data = New-Object 'object[,]' 1,18
    for ($col = 0; $col -lt 18; $col++) {
        $data[0,$col] = $rowValues[$col]
    }
    $ws.Range("O$rowIdx`:AF$rowIdx").Value = $data
    $rowIdx++
}

# --- Selection: move from Y9 to AC30, drop the old topLeftCell=Q1 scroll --
$ws.Range("AC30").Select()
